# Update latest output (run 88)
$wb = $excel.ActiveWorkbook

# ---- Sheet "Schedule" ----
$ws1 = $wb.Worksheets.Item("Schedule")

$ws1.Range("A2").Value = 46040.25
$ws1.Range("C2").Value = 14.5
$ws1.Range("D2").Value = 54.81
$ws1.Range("E2").Value = 143.6525415
$ws1.Range("F2").Value = 2.620918472906404

$ws1.Range("A3").Value = 46040.9375
$ws1.Range("C3").Value = 4.5
$ws1.Range("D3").Value = 17.01
$ws1.Range("E3").Value = 398.40378825
$ws1.Range("F3").Value = 23.42173946208113

$ws1.Range("E4").Value = 69.74660550000003
$ws1.Range("F4").Value = 2.050164770723105

# ---- Sheet "Detailed" ----
$ws2 = $wb.Worksheets.Item("Detailed")

$ws2.Range("E14").Value = "ON"

$ws2.Range("B40").Value = 36.2
$ws2.Range("B42").Value = 46.53455
$ws2.Range("B43").Value = 62.21582
$ws2.Range("C43").Value = "historical"
$ws2.Range("B44").Value = 62.97041
$ws2.Range("C44").Value = "historical"
$ws2.Range("B46").Value = 57.06007
$ws2.Range("E46").Value = "OFF"
$ws2.Range("B47").Value = 56.98
$ws2.Range("B48").Value = 55.57436
$ws2.Range("B49").Value = 49.23901
$ws2.Range("B50").Value = 49.34509
$ws2.Range("B51").Value = 55.94746
$ws2.Range("B52").Value = 40.54
$ws2.Range("B53").Value = 35.86
$ws2.Range("B54").Value = 29.25339
$ws2.Range("B55").Value = 35.87996
$ws2.Range("B56").Value = 36.2
$ws2.Range("B57").Value = 45.21307
$ws2.Range("B59").Value = 58.71596
$ws2.Range("B60").Value = 58.95696
$ws2.Range("B61").Value = 65
$ws2.Range("B64").Value = 36.06
$ws2.Range("B66").Value = 24.92768
$ws2.Range("B67").Value = 22.07
$ws2.Range("B68").Value = 0.51
$ws2.Range("B69").Value = -5.01
$ws2.Range("B70").Value = -5.50985
$ws2.Range("B71").Value = -5.79092
$ws2.Range("B72").Value = -4.83724
$ws2.Range("B73").Value = -5.50985
$ws2.Range("B75").Value = -4.85271
$ws2.Range("B76").Value = -5.01
$ws2.Range("B77").Value = -2.5711
$ws2.Range("B78").Value = 0
$ws2.Range("B79").Value = -0.91982
$ws2.Range("B80").Value = -2.57161
$ws2.Range("B81").Value = 0.00025
$ws2.Range("B82").Value = 0.00048
$ws2.Range("B83").Value = -5.25561
$ws2.Range("B84").Value = -2.53803
$ws2.Range("B85").Value = 0.01346
$ws2.Range("B86").Value = 12.05992
$ws2.Range("B87").Value = 52.67227
$ws2.Range("B88").Value = 57.46321
$ws2.Range("B89").Value = 62.05077
$ws2.Range("B90").Value = 73.20007
$ws2.Range("B92").Value = 69.62199
$ws2.Range("B94").Value = 62.98392
$ws2.Range("B95").Value = 58.98682
$ws2.Range("B96").Value = 58.15428
$ws2.Range("B97").Value = 61.96632
